$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty spacer column F (width 50.5); this shifts column G into F.
$ws.Columns("F").Delete()

# Swap columns B and C (code column now precedes the name column).
$ws.Columns("C").Cut()
$ws.Columns("B").Insert()

# Swap columns D and E (code column now precedes the name column).
$ws.Columns("E").Cut()
$ws.Columns("D").Insert()

# Update the active selection to match the saved view state.
$ws.Range("C20").Select()
